$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.215181708335876
$ws.Range("B1").Value = 2.440188646316528
$ws.Range("C1").Value = 7.176395416259766
$ws.Range("D1").Value = 2.255922079086304
$ws.Range("E1").Value = 1.164646625518799
